# Updated cryptos list (price/volume refresh) per upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.335.20"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.841.04"
$ws.Range("E3").Value = "  -0.52%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.14"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$ws.Range("E6").Value = "  -1.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07417"
$ws.Range("E8").Value = "  -2.09%  "

# Row 9
$ws.Range("E9").Value = "  -0.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.80"
$ws.Range("E10").Value = "  +1.27%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07732"
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "1.841.83"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.980"
$ws.Range("E13").Value = "  -0.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6771"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001019"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.98"
$ws.Range("E16").Value = "  -1.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.246"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18
$ws.Range("D18").Value = "29.372.03"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.96"
$ws.Range("E19").Value = "  -0.50%  "

# Row 20
$ws.Range("E20").Value = "  -0.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.416"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.05"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.468"
$ws.Range("E25").Value = "  +0.15%  "

# Row 26
$ws.Range("E26").Value = "  -3.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.42"
$ws.Range("E27").Value = "  -1.42%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.446"
$ws.Range("E28").Value = "  +2.08%  "

# Row 29
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06438"
$ws.Range("E29").Value = "  +13.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.486"
$ws.Range("E30").Value = "  +0.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.066"
$ws.Range("E31").Value = "  -1.72%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.068"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.835"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("E34").Value = "  -1.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6913"
$ws.Range("E35").Value = "  -1.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.559"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01857"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.822"
$ws.Range("E38").Value = "  +3.48%  "

# Row 39
$ws.Range("D39").Value = "1.241.96"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.737"
$ws.Range("E40").Value = "  +3.62%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9307"
$ws.Range("E41").Value = "  +2.96%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "2.026.28"
$ws.Range("E43").Value = "  +0.75%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.71"
$ws.Range("E44").Value = "  -1.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.66"
$ws.Range("E45").Value = "  -0.43%  "

# Row 46
$ws.Range("E46").Value = "  +2.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.054"
$ws.Range("E47").Value = "  -1.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.715"
$ws.Range("E48").Value = "  +2.28%  "

# Row 49
$ws.Range("E49").Value = "  -1.54%  "

# Row 50
$ws.Range("E50").Value = "  -0.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3888"
$ws.Range("E51").Value = "  -1.94%  "

